# Update the "想去人数" (want-to-go count) figures in column F
# for both the "展览" and "全部类型" worksheets.
# Rows and their new values, as scraped from the upstream data source refresh.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 263
    5  = 835
    6  = 2
    7  = 292
    8  = 7374
    9  = 70
    12 = 100
    13 = 3
    14 = 43
    18 = 247
    19 = 666
    20 = 17
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
